$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of CRM lead data appended below the existing rows
$ws.Range("A4").Value = "FIS"
$ws.Range("B4").Value = "Vibhuti"
$ws.Range("C4").Value = "Mishra"

# Leading apostrophe forces text storage so the phone number keeps the
# same "quote-prefixed text" style used by the other PhoneNumber cells.
$ws.Range("D4").Value = "'8979466578"

$ws.Range("E4").Value = "vibhuti.mishra@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:vibhuti.mishra@gmail.com")
$ws.Range("E4").Style = "Hyperlink"
